# Generate Report for Handoff
# Rename source file b742e51e-0df9-44be-a16e-1a022713b4da.md -> 7902faeb-a6e2-40b0-838f-2b1dcb490321.md
# everywhere it is referenced, refresh the handoff/generate timestamps, and clear
# the (now stale) handback target/file/datetime columns on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$oldName = "b742e51e-0df9-44be-a16e-1a022713b4da"
$newName = "7902faeb-a6e2-40b0-838f-2b1dcb490321"

$oldXliffHash = "25b5561527b3d21c1c8e4884f0d0b37954a58ac1"
$newXliffHash = "c50cd22493fc31116525de9d2a0a81778cf89756"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newName.md"
$wsOverview.Range("B2").Value = "e2e\$newName.md"
$wsOverview.Range("G2").Value = "2016-08-15 20:55:59"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newName.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newName.md"
$wsZhCn.Range("G2").Value = "$newName.$newXliffHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-15 20:55:53"
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Clear the stale "Latest Target File" / "Latest Handback File" cells
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"

# Keep the A2 hyperlink (pointing at the source file) but drop the I2 one
# (the handback file hyperlink is no longer valid) -- update A2's display text.
$zhHyperlinks = @()
foreach ($hl in $wsZhCn.Hyperlinks) { $zhHyperlinks += $hl }
foreach ($hl in $zhHyperlinks) {
    if ($hl.Range.Address() -eq '$I$2') {
        $hl.Delete()
    } elseif ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newName.md"
    }
}

$wsZhCn.Columns.Item(9).ColumnWidth = 17.83
$wsZhCn.Columns.Item(10).ColumnWidth = 20.83

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newName.md"
$wsDeDe.Range("G2").Value = "$newName.$newXliffHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-15 20:55:59"
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"

$deHyperlinks = @()
foreach ($hl in $wsDeDe.Hyperlinks) { $deHyperlinks += $hl }
foreach ($hl in $deHyperlinks) {
    if ($hl.Range.Address() -eq '$I$2') {
        $hl.Delete()
    } elseif ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newName.md"
    }
}

$wsDeDe.Columns.Item(9).ColumnWidth = 17.83
$wsDeDe.Columns.Item(10).ColumnWidth = 20.83
